# Gantt-Diagramm.xlsx: "Pflichtenheft zusammengefasst und Gantt erweitert"
#
# 1. Advance the Gantt chart's display-week selector (Anzeigewoche, linked
#    to cell E4 via the scrollbar control) from week 8 to week 12 - this
#    shifts the whole displayed 10-week window by 4 weeks (28 days); every
#    date in rows 4/5 is formula-driven off this cell and recalculates
#    automatically.
# 2. Update task progress: "Arbeitspakete definieren" (row 12) -> 100%,
#    "7, Organigramm erstellen" (row 26) -> reset to 0%.
# 3. Scroll/re-select the visible window further down the task list.
# 4. Slightly reduce the print scale.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projektplan")

# 1. Anzeigewoche: 8 -> 12 (drives the whole Gantt header via formulas)
$ws.Range("E4").Value = 12

# Keep the scrollbar's own stored position in sync with the linked cell,
# best-effort (the control mirrors $E$4).
$ws.Shapes.Item(1).ControlFormat.Value = 12

# 2. Task progress updates
$ws.Range("D12").Value = 1    # Julian & Gabriel: Arbeitspakete definieren -> 100%
$ws.Range("D26").Value = 0    # Gabriel: 7, Organigramm erstellen -> 0%

# 3. View: scroll the frozen pane further down and move the selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$ws.Range("D2").Select()

# 4. Page setup: tighten the print scale slightly (60% -> 59%)
$ws.PageSetup.Zoom = 59
$ws.PageSetup.FitToPagesTall = $false
